$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.070.39"
$ws.Range("E2").Value = "  +0.30%  "

$ws.Range("D3").Value = "2.842.57"
$ws.Range("E3").Value = "  +2.56%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "363.40"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +7.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "112.87"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.15%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.571"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +4.80%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.604"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +5.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.52"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.27%  "

$ws.Range("E11").Value = "  +0.82%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.24"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.51%  "

$ws.Range("E13").Value = "  +1.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.79"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.92%  "

$ws.Range("D15").Value = "3.289.04"
$ws.Range("E15").Value = "  +2.59%  "

$ws.Range("D16").Value = "2.860.09"
$ws.Range("E16").Value = "  +3.00%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.925"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +5.71%  "

$ws.Range("D18").Value = "52.020.21"
$ws.Range("E18").Value = "  +0.50%  "

$ws.Range("E19").Value = "  +9.01%  "

$ws.Range("E20").Value = "  -0.46%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.62"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.42%  "

$ws.Range("E22").Value = "  +2.78%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.48"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.69"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.74%  "

$ws.Range("E25").Value = "  +4.71%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.17"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.18%  "

$ws.Range("E27").Value = "  +0.04%  "

$ws.Range("E28").Value = "  +2.54%  "

$ws.Range("E29").Value = "  +1.34%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0488"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +30.45%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "53.48"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +6.74%  "

$ws.Range("E32").Value = "  +0.26%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.64"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.40%  "

$ws.Range("E34").Value = "  +3.45%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.53"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +12.38%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0848"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.31%  "

$ws.Range("E37").Value = "  -0.12%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.30"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.96%  "

$ws.Range("E39").Value = "  -0.96%  "

$ws.Range("E40").Value = "  -1.70%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.73"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.53%  "

$ws.Range("E42").Value = "  +1.93%  "

$ws.Range("E43").Value = "  -4.19%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "126.34"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.41%  "

$ws.Range("E45").Value = "  -2.28%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.41"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.65%  "

$ws.Range("D47").Value = "2.115.74"
$ws.Range("E47").Value = "  +2.19%  "

$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.95"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +8.16%  "

$ws.Range("B50").Value = "SEI"
$ws.Range("C50").Value = "https://coinranking.com/coin/8nxCqs-uj+sei-sei"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.990"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +13.67%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "61.89"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +4.85%  "
